$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: set C3 = 0 (B3 stays 0)
$ws.Range("C3").Value = 0

# Rows 4-20: update B (measured Vi) and C (measured Vo) values
$data = @(
    @(4, 0.106, 0.106),
    @(5, 0.205, 0.205),
    @(6, 0.305, 0.304),
    @(7, 0.355, 0.351),
    @(8, 0.404, 0.393),
    @(9, 0.455, 0.427),
    @(10, 0.504, 0.451),
    @(11, 0.554, 0.469),
    @(12, 0.603, 0.483),
    @(13, 0.654, 0.495),
    @(14, 0.704, 0.505),
    @(15, 0.753, 0.513),
    @(16, 0.803, 0.52),
    @(17, 0.902, 0.532),
    @(18, 1.002, 0.541),
    @(19, 1.499, 0.573),
    @(20, 1.995, 0.592)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
}

$ws.Range("G4").Select()
